# The source data (array of feature columns B:J, one row per centroid) was
# re-indexed using fancy/array indexing (reshape with double brackets [[ ]]),
# which reorders the centroid rows. Column A (the centroid id) and the header
# row stay put; only the feature values in columns B:J for rows 2-11 move.
#
# New row order expressed as the 0-based original-row index that now supplies
# the B:J values for each destination row (rows 2..11 <-> ids 0..9):
#   dest id 0 <- src id 2
#   dest id 1 <- src id 1
#   dest id 2 <- src id 3
#   dest id 3 <- src id 7
#   dest id 4 <- src id 4
#   dest id 5 <- src id 5
#   dest id 6 <- src id 6
#   dest id 7 <- src id 9
#   dest id 8 <- src id 8
#   dest id 9 <- src id 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:J content for each row (row 2 through row 11), already reordered.
$newRows = @(
    @(-0, 1, 1, 0, 0, 0, -0, 0.5604, 44.2716),
    @(1, 0, 1, -0, 0, 0, -0, 0.5425, 42.8575),
    @(0, 1, 0, -0, 0, 1, -0, 0.5513, 43.5527),
    @(1, 0, 0, -0, 0, -0, 1, 0.5661, 44.72190000000001),
    @(0, 1, 0, 1, 0, -0, 0, 0.5713, 45.1327),
    @(1, 0, 0, 1, 0, -0, -0, 0.5537, 43.7423),
    @(0, 1, -0, -0, -0, -0, 1, 0.5339, 42.1781),
    @(1, 0, 0, -0, 0, 1, -0, 0.5456, 43.1024),
    @(1, 0, 0, -0, 1, 0, -0, 0.5737, 45.3223),
    @(0, 1, 0, -0, 1, 0, -0, 0.5845, 46.1755)
)

$startRow = 2
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowValues = $newRows[$i]
    $excelRow = $startRow + $i
    for ($col = 0; $col -lt $rowValues.Length; $col++) {
        # Columns B..J are Excel columns 2..10
        $ws.Cells.Item($excelRow, $col + 2).Value = $rowValues[$col]
    }
}
